$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 3; $row -le 15; $row++) {
    $ws.Cells.Item($row, 5).Formula = '=AVERAGE(B' + $row + ':D' + $row + ')'
    $ws.Cells.Item($row, 6).Formula = '=STDEV.S(B' + $row + ':D' + $row + ')'
    $ws.Cells.Item($row, 7).Formula = '=CONCATENATE("{",A' + $row + ',",",B' + $row + ',"},{",A' + $row + ',",",C' + $row + ',"},{",A' + $row + ',",",D' + $row + ',"},")'
}
